$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ----------------------------------------------------------------------
# Formatting is propagated by copy/paste-special (formats only) from
# existing cells that already carry the exact xf this edit needs, so we
# never have to hand-roll style indices. Order matters: G3's current
# format (border-only, no extra fill) is harvested first, before G3 is
# itself reformatted a few lines down.
# ----------------------------------------------------------------------

# "thin border, no fill" look currently on G3:G6 -> needed on C/E/G 27:29
$ws.Range("G3").Copy()
$ws.Range("C27:C29").PasteSpecial(-4122)
$ws.Range("E26:E29").PasteSpecial(-4122)
$ws.Range("G27:G29").PasteSpecial(-4122)

# plain border + fill flag look currently on C8 -> needed on G2, D26:D29
$ws.Range("C8").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("D26:D29").PasteSpecial(-4122)

# header look currently on A1 -> needed on G1
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# look currently on F2 -> needed on F26:F29
$ws.Range("F2").Copy()
$ws.Range("F26:F29").PasteSpecial(-4122)

# plain thin-border look currently on A2 -> needed on G3:G26 and the
# remaining new-row cells
$ws.Range("A2").Copy()
$ws.Range("G3:G26").PasteSpecial(-4122)
$ws.Range("A26:C26").PasteSpecial(-4122)
$ws.Range("A27:B29").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# Value updates
# ----------------------------------------------------------------------

# icon paths that were re-exported with spaced-out SVG path commands
$ws.Range("G4").Value = 'path("M 0 8.833 V 0.5 A 0.5 0.5 0 0 1 0.5 0 H 8.833 A 0.5 0.5 0 0 1 9.333 0.5 V 8.833 A 0.5 0.5 0 0 1 8.833 9.333 H 0.5 A 0.5 0.5 0 0 1 0 8.833 Z M 3.468 11.617 V 10.426 H 9.853 A 0.55 0.55 0 0 0 10.425 9.853 V 3.468 H 11.617 V 11.618 H 3.468 Z M 14 14 H 5.851 V 12.809 H 12.236 A 0.55 0.55 0 0 0 12.808 12.236 V 5.851 H 14 V 14 Z")'
$ws.Range("G5").Value = 'path("M 0.215 10.83 C -0.374 9.51 0.288 7.969 1.686 7.404 L 4.42 6.301 L 6.788 11.608 L 6.518 11.718 L 7.508 13.935 L 4.87 15 L 3.873 12.767 C 2.499 13.216 1.004 12.597 0.443 11.339 L 0.215 10.829 Z M 11.132 0.012 C 11.302 0.052 11.442 0.155 11.509 0.306 L 15.95 10.259 C 15.977 10.319 16.004 10.379 16 10.453 C 16.005 10.719 15.765 10.942 15.485 10.947 L 7.135 11.476 L 4.772 6.15 L 10.664 0.148 A 0.507 0.507 0 0 1 11.132 0.013 Z")'
$ws.Range("G7").Value = 'path("M 7 0 C 10.862 0 14 3.138 14 7 S 10.862 14 7 14 S 0 10.862 0 7 S 3.138 0 7 0 Z M 8.918 9.67 L 8.556 7.518 L 10.111 5.988 L 7.959 5.678 L 7 3.732 L 6.04 5.676 L 3.889 5.987 L 5.444 7.517 L 5.082 9.669 L 7 8.659 L 8.918 9.67 Z")'

# four new icon test-case rows
$ws.Range("A26").Value = 'TC_025'
$ws.Range("B26").Value = 'Verify the help icon and Css Values'
$ws.Range("C26").Value = 'Help'
$ws.Range("D26").Value = '14px'
$ws.Range("E26").Value = 'var(--cool-grey)'
$ws.Range("F26").Value = 'rgba(21, 78, 197, 1)'
$ws.Range("G26").Value = 'path("m 10 0 a 10 10 0 1 1 10 20.002 a 10 10 0 0 1 10 0 z m 9.098 14.777 h 11.028 v 12.925 h 9.098 v 14.777 z m 12.303 5.922 a 4.015 4.015 0 0 0 9.906 5.237 a 3.36 3.36 0 0 0 8.039 5.736 c 7.229 6.255 6.804 7.116 6.762 8.319 h 8.63 c 8.632 7.953 8.74 7.597 8.941 7.292 a 1.15 1.15 0 0 1 10 6.793 c 10.498 6.793 10.84 6.928 11.027 7.198 c 11.216 7.461 11.314 7.778 11.307 8.1 c 11.309 8.38 11.222 8.653 11.057 8.879 c 10.954 9.023 10.827 9.149 10.684 9.253 l 10.218 9.626 a 2.658 2.658 0 0 0 9.346 10.591 c 9.199 11.076 9.126 11.579 9.128 12.085 h 10.872 a 2.8 2.8 0 0 1 10.965 11.339 c 11.041 11.056 11.204 10.805 11.432 10.622 l 11.899 10.279 c 12.247 10.033 12.562 9.739 12.833 9.409 c 13.113 8.997 13.255 8.505 13.237 8.009 a 2.364 2.364 0 0 0 12.303 5.922 z")'

$ws.Range("A27").Value = 'TC_026'
$ws.Range("B27").Value = 'Verify the tender icon icon and Css Values'
$ws.Range("C27").Value = 'Tender'
$ws.Range("D27").Value = '14px'
$ws.Range("E27").Value = 'var(--cool-grey)'
$ws.Range("F27").Value = 'rgba(21, 78, 197, 1)'
$ws.Range("G27").Value = 'path("M 4 20 L 5 11 H 0 L 2 0 H 13 L 9 7 H 14 Z")'

$ws.Range("A28").Value = 'TC_027'
$ws.Range("B28").Value = 'Verify the Notification icon and Css Values'
$ws.Range("C28").Value = 'Notification'
$ws.Range("D28").Value = '14px'
$ws.Range("E28").Value = 'var(--cool-grey)'
$ws.Range("F28").Value = 'rgba(21, 78, 197, 1)'
$ws.Range("G28").Value = 'path("M 7.313 2.388 V 1.636 C 7.312 0.733 8.07 0 9 0 C 9.932 0 10.687 0.73 10.687 1.636 V 2.388 C 13.6 3.115 15.75 5.678 15.75 8.728 V 13.091 C 15.75 13.494 16.085 13.818 16.499 13.818 C 17.325 13.818 18 14.468 18 15.273 V 16.365 A 0.371 0.371 0 0 1 17.624 16.727 H 0.376 A 0.368 0.368 0 0 1 0 16.365 V 15.273 C 0 14.472 0.672 13.818 1.501 13.818 A 0.74 0.74 0 0 0 2.25 13.091 V 8.727 C 2.25 5.677 4.402 3.114 7.313 2.387 Z M 6.375 17.455 H 11.625 C 11.625 18.86 10.45 20 9 20 S 6.375 18.86 6.375 17.455 Z")'

$ws.Range("A29").Value = 'TC_028'
$ws.Range("B29").Value = 'Verify the Profile icon and Css Values'
$ws.Range("C29").Value = 'Profile'
$ws.Range("D29").Value = '14px'
$ws.Range("E29").Value = 'var(--cool-grey)'
$ws.Range("F29").Value = 'rgba(21, 78, 197, 1)'
$ws.Range("G29").Value = 'path("M 7.313 2.388 V 1.636 C 7.312 0.733 8.07 0 9 0 C 9.932 0 10.687 0.73 10.687 1.636 V 2.388 C 13.6 3.115 15.75 5.678 15.75 8.728 V 13.091 C 15.75 13.494 16.085 13.818 16.499 13.818 C 17.325 13.818 18 14.468 18 15.273 V 16.365 A 0.371 0.371 0 0 1 17.624 16.727 H 0.376 A 0.368 0.368 0 0 1 0 16.365 V 15.273 C 0 14.472 0.672 13.818 1.501 13.818 A 0.74 0.74 0 0 0 2.25 13.091 V 8.727 C 2.25 5.677 4.402 3.114 7.313 2.387 Z M 6.375 17.455 H 11.625 C 11.625 18.86 10.45 20 9 20 S 6.375 18.86 6.375 17.455 Z")'

$ws.Range("D17").Select()
